$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set all target cells to text format first so numeric-looking strings
# (e.g. "1.000", "25.064.75") are preserved exactly as text, matching the
# source data which stores these as literal strings, not numbers/dates.
$targetCells = @(
  "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D7", "E7", "B8", "C8", "D8", "E8", "B9", "C9",
  "D9", "E9", "B10", "C10", "D10", "E10", "B11", "C11", "D11", "E11", "B12", "C12", "D12", "E12",
  "B13", "C13", "D13", "E13", "B14", "C14", "D14", "E14", "B15", "C15", "D15", "E15", "B16", "C16",
  "D16", "E16", "B17", "C17", "E17", "B18", "C18", "D18", "E18", "B19", "C19", "D19", "E19", "B20",
  "C20", "D20", "E20", "B21", "C21", "D21", "E21", "B22", "C22", "D22", "E22", "B23", "C23", "D23",
  "E23", "B24", "C24", "D24", "E24", "B25", "C25", "D25", "E25", "B26", "C26", "D26", "E26", "B27",
  "C27", "D27", "E27", "B28", "C28", "D28", "E28", "B29", "C29", "D29", "E29", "B30", "C30", "D30",
  "E30", "B31", "C31", "D31", "E31", "B32", "C32", "D32", "E32", "B33", "C33", "D33", "E33", "B34",
  "C34", "D34", "E34", "B35", "C35", "D35", "E35", "B36", "C36", "D36", "E36", "B37", "C37", "D37",
  "E37", "B38", "C38", "D38", "E38", "B39", "C39", "D39", "E39", "B40", "C40", "D40", "E40", "B41",
  "C41", "D41", "E41", "B42", "C42", "D42", "E42", "B43", "C43", "D43", "E43", "B44", "C44", "D44",
  "E44", "B45", "C45", "D45", "E45", "B46", "C46", "D46", "E46", "B47", "C47", "D47", "E47", "B48",
  "C48", "D48", "E48", "B49", "C49", "D49", "E49", "B50", "C50", "D50", "E50", "B51", "C51", "D51",
  "E51"
)
foreach ($cellRef in $targetCells) {
  $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "25.064.75"
$ws.Range("E2").Value = "  -3.61%  "
$ws.Range("D3").Value = "1.649.59"
$ws.Range("E3").Value = "  -5.51%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "237.81"
$ws.Range("E5").Value = "  -4.50%  "
$ws.Range("D7").Value = "0.4795"
$ws.Range("E7").Value = "  -6.77%  "
$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").Value = "39.63"
$ws.Range("E8").Value = "  -2.70%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "0.2616"
$ws.Range("E9").Value = "  -5.17%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.06012"
$ws.Range("E10").Value = "  -3.05%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "0.07170"
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.653.05"
$ws.Range("E12").Value = "  -5.32%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "14.83"
$ws.Range("E13").Value = "  -2.40%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "0.6237"
$ws.Range("E14").Value = "  -3.88%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "4.597"
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "73.35"
$ws.Range("E16").Value = "  -5.80%  "
$ws.Range("B17").Value = "Dai"
$ws.Range("C17").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("B18").Value = "BinanceUSD"
$ws.Range("C18").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D18").Value = "0.9999"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "25.053.93"
$ws.Range("E19").Value = "  -3.76%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "11.48"
$ws.Range("E20").Value = "  -3.19%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.000006622"
$ws.Range("E21").Value = "  -2.97%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "4.483"
$ws.Range("E22").Value = "  +4.33%  "
$ws.Range("B23").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C23").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D23").Value = "1.861.20"
$ws.Range("E23").Value = "  -5.28%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "8.603"
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").Value = "5.293"
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "132.06"
$ws.Range("E26").Value = "  -2.40%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "14.98"
$ws.Range("E27").Value = "  -2.18%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "1.395"
$ws.Range("E28").Value = "  -7.33%  "
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "103.40"
$ws.Range("E29").Value = "  -2.60%  "
$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").Value = "1.683"
$ws.Range("E30").Value = "  -5.66%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "3.773"
$ws.Range("E31").Value = "  -4.67%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "0.07908"
$ws.Range("E32").Value = "  -4.27%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "3.582"
$ws.Range("E33").Value = "  -2.64%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.04595"
$ws.Range("E34").Value = "  -1.70%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.599"
$ws.Range("E35").Value = "  -2.13%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "0.9425"
$ws.Range("E36").Value = "  -5.67%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "0.5776"
$ws.Range("E37").Value = "  -7.37%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "2.612"
$ws.Range("E38").Value = "  -4.55%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01555"
$ws.Range("E39").Value = "  -3.11%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").Value = "1.001"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "0.8206"
$ws.Range("E41").Value = "  +8.45%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "1.826"
$ws.Range("E42").Value = "  -5.45%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "98.63"
$ws.Range("E43").Value = "  -1.68%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.3727"
$ws.Range("E44").Value = "  -4.05%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "4.810"
$ws.Range("E45").Value = "  -4.22%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "0.1143"
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "6.109"
$ws.Range("E47").Value = "  -3.71%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.05188"
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "29.79"
$ws.Range("E49").Value = "  -2.93%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "50.90"
$ws.Range("E50").Value = "  -8.21%  "
$ws.Range("B51").Value = "TrueUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
$ws.Range("D51").Value = "1.002"
$ws.Range("E51").Value = "  -0.17%  "
